{"js": "// Change the time-format placeholder from hyphen-separated ('hh-mm-ss')\n// to colon-separated ('hh:mm:ss') everywhere it appears in the document\n// (Starting_time and Ending_time bullet points).\nconst body = context.document.body;\nconst results = body.search(\"hh-mm-ss\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"hh:mm:ss\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Change the time-format placeholder from hyphen-separated ('hh-mm-ss')\n# to colon-separated ('hh:mm:ss') everywhere it appears in the document\n# (Starting_time and Ending_time bullet points).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find.Execute(\"hh-mm-ss\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"hh:mm:ss\", $wdReplaceAll)\n"}
